# DG: updated Implementation Undo/Redo Steps 1 - 4
#
# Renames the ":AddressBookParser" / "AddressBook" / "VersionedAddressBook" /
# "ReadOnlyAddressBook" class names in the sequence diagram to the
# "FinanceTrackerParser" / "Finance Tracker" / "VersionedFinanceTracker" /
# "ReadOnlyFinanceTracker" equivalents, repositioning/resizing the boxes that
# grew to fit the new (longer) text.

function EMUToPt($emu) {
    # Shape.Left/Top/Width/Height are expressed in points; the interop layer
    # truncates the float -> EMU conversion, so nudge up slightly to land on
    # the exact EMU value we want.
    return ($emu / 12700.0) + 0.00003
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# Shape "Rectangle 62" (id 16) - ":Address" / "BookParser" -> ":FinanceTrackerParser"
# ---------------------------------------------------------------------
$sh = $s.Shapes.Item(6)
$sh.Left = EMUToPt 3307635
$sh.Top = EMUToPt 446716
$sh.Width = EMUToPt 1492965
$sh.Height = EMUToPt 467684

$tr = $sh.TextFrame.TextRange
$tr.Text = ":" + "FinanceTrackerParser"
$run2 = $tr.Characters(2, 20)
$run2.Font.Bold = $run2.Font.Bold

# ---------------------------------------------------------------------
# Shape "TextBox 78" (id 79) - "undoAddressBook()" -> "undoFinanceTracker()"
# (now wraps onto two lines: "undoFinance" / "Tracker()")
# ---------------------------------------------------------------------
$sh = $s.Shapes.Item(19)
$sh.Left = EMUToPt 5386369
$sh.Top = EMUToPt 2599418
$sh.Width = EMUToPt 1298078
$sh.Height = EMUToPt 369332

$tr = $sh.TextFrame.TextRange
$tr.Text = "undoFinance" + [char]13 + "Tracker()"
$run2 = $tr.Characters(5, 7)
$run2.Font.Color.RGB = 0xA03070
$run3 = $tr.Characters(13, 9)
$run3.Font.Color.RGB = 0xA03070

# ---------------------------------------------------------------------
# Shape "Rectangle 62" (id 84) - ":VersionedAddressBook" -> ":VersionedFinanceTracker"
# ---------------------------------------------------------------------
$sh = $s.Shapes.Item(23)
$sh.Left = EMUToPt 7497155
$sh.Top = EMUToPt 2505472
$sh.Width = EMUToPt 2324727
$sh.Height = EMUToPt 398561

$tr = $sh.TextFrame.TextRange
$full = $tr.Text
$run2 = $tr.Characters(2, $full.Length - 1)
$run2.Text = "VersionedFinanceTracker"

# ---------------------------------------------------------------------
# Shape "Rectangle 62" (id 40) - ": Model" box just moves (no text change)
# ---------------------------------------------------------------------
$sh = $s.Shapes.Item(26)
$sh.Left = EMUToPt 6549220
$sh.Top = EMUToPt 2291613

# ---------------------------------------------------------------------
# Shape "TextBox 87" (id 88) - "resetData(ReadOnlyAddressBook)" -> "resetData(ReadOnlyFinanceTracker)"
# ---------------------------------------------------------------------
$sh = $s.Shapes.Item(35)
$sh.Left = EMUToPt 8769529
$sh.Top = EMUToPt 3267337
$sh.Width = EMUToPt 2321759
$sh.Height = EMUToPt 184666

$tr = $sh.TextFrame.TextRange
$run2 = $tr.Characters(11, 19)
$run2.Text = "ReadOnlyFinanceTracker"
